$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PO List")

# --- Updated cell values on "PO List" (row 3 .. 29) ---
$ws.Range("N3").Value = 15
$ws.Range("R3").Value = 16
$ws.Range("R4").Value = 4
$ws.Range("Y4").Value = 3
$ws.Range("R5").Value = 2
$ws.Range("R6").Value = 18
$ws.Range("N7").Value = 4
$ws.Range("R7").Value = 7
$ws.Range("R8").Value = 8
$ws.Range("R9").Value = 10
$ws.Range("Y9").Value = 7
$ws.Range("N10").Value = 14
$ws.Range("R10").Value = 14
$ws.Range("N11").Value = 12
$ws.Range("R11").Value = 6
$ws.Range("Y11").Value = 4
$ws.Range("K12").Value = 9
$ws.Range("L12").Value = 9
$ws.Range("N12").Value = 10
$ws.Range("O12").Value = 5
$ws.Range("P12").Value = 5
$ws.Range("Q12").Value = 44981
$ws.Range("R12").Value = 3
$ws.Range("S12").Value = 4
$ws.Range("T12").Value = 4
$ws.Range("U12").Value = 44981
$ws.Range("V12").Value = 8
$ws.Range("W12").Value = 8
$ws.Range("X12").Value = 44981
$ws.Range("Y12").Value = 2
$ws.Range("N13").Value = 9
$ws.Range("F14").Value = 2
$ws.Range("R15").Value = 15
$ws.Range("K16").Value = 9
$ws.Range("L16").Value = 9
$ws.Range("M16").Value = 45007
$ws.Range("N16").Value = 1
$ws.Range("R16").Value = 5
$ws.Range("V16").Value = 11
$ws.Range("W16").Value = 11
$ws.Range("X16").Value = 45007
$ws.Range("Y16").Value = 1
$ws.Range("R17").Value = 17
$ws.Range("N18").Value = 2
$ws.Range("R18").Value = 12
$ws.Range("N19").Value = 7
$ws.Range("Y19").Value = 6
$ws.Range("N20").Value = 4
$ws.Range("R20").Value = 19
$ws.Range("R21").Value = 21
$ws.Range("N24").Value = 11
$ws.Range("R24").Value = 20
$ws.Range("N25").Value = 8
$ws.Range("R25").Value = 9
$ws.Range("R26").Value = 13
$ws.Range("N27").Value = 6
$ws.Range("O27").Value = 3
$ws.Range("P27").Value = 3
$ws.Range("Q27").Value = 45005
$ws.Range("R27").Value = 1
$ws.Range("S27").Value = 1
$ws.Range("T27").Value = 1
$ws.Range("U27").Value = 45005
$ws.Range("N28").Value = 3
$ws.Range("Y28").Value = 5
$ws.Range("N29").Value = 13
$ws.Range("R29").Value = 11

# --- View state: make "PO List" the active/selected tab ---
$ws.Activate()
$ws.Range("B3:Z29").Select()
